$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The timesheet sheet is being turned into a task-hours summary:
#   - the old per-punch header row (Employee/Date/Clock In/Clock Out/
#     Hours Worked) goes away
#   - row 2 becomes a spacer row whose last column (G) reports the
#     grand total instead of the old "Total Hours:" / 0 pair
#   - a small "Task Name" / "Total Hours" mini-table header is added
#     a couple of rows below
# ------------------------------------------------------------------

# Grab the formatting (bold, no centering) that the original D3/E3
# placeholder cells used, *before* touching anything, so every new
# bold cell below reuses that existing style instead of Excel minting
# new ones.
$ws.Range("D3").Copy()

# New grand-total cell, column G of the spacer row.
$ws.Range("G2").Value = "Overall Total Hours: 0.0"
$ws.Range("G2").PasteSpecial(-4122)

# Blank styled placeholder a few rows down (row 5) ahead of the new
# mini-table.
$ws.Range("A5").PasteSpecial(-4122)

# "Task Name" / "Total Hours" column headers for the new table.
$ws.Range("A6").Value = "Task Name"
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = "Total Hours"
$ws.Range("B6").PasteSpecial(-4122)

# Drop the old employee clock-in/out header row entirely.
$ws.Range("A1:E1").Clear()

# The old "Total Hours:" label + value on row 2 are replaced by the
# new G2 summary above, so clear them out.
$ws.Range("D2:E2").Clear()

# The old styled placeholder cells on row 3 are no longer needed now
# that their format has been copied forward.
$ws.Range("D3:E3").Clear()

# Two more columns (F, G) are now part of the sheet's used range;
# give them the same "13"-wide layout as the existing columns.
$ws.Columns.Item(6).ColumnWidth = 12.2
$ws.Columns.Item(7).ColumnWidth = 12.2

Write-Output "done"
